$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1065.2572
$ws.Range("I33").Value = 302.9
$ws.Range("J33").Value = 5639.4
$ws.Range("K33").Value = 302.9
$ws.Range("L33").Value = 5639.4
$ws.Range("M33").Value = -73.89999999999998
$ws.Range("N33").Value = -6097.4

$ws.Range("H98").Value = 2260.3076
$ws.Range("I98").Value = 1217.6364
$ws.Range("J98").Value = 7995
$ws.Range("K98").Value = 1217.6364
$ws.Range("L98").Value = 7995
$ws.Range("M98").Value = 280.3635999999999
$ws.Range("N98").Value = -10991

$ws.Range("H122").Value = 2260.3076
$ws.Range("I122").Value = 1217.6364
$ws.Range("J122").Value = 7995
$ws.Range("K122").Value = 3652.9092
$ws.Range("L122").Value = 23985
$ws.Range("M122").Value = -1202.9092
$ws.Range("N122").Value = -28885

$ws.Range("H129").Value = 2624.4465
$ws.Range("J129").Value = 930.907
$ws.Range("L129").Value = 2792.721
$ws.Range("N129").Value = -12792.721

$ws.Range("H141").Value = 1868.3489
$ws.Range("I141").Value = 1771.8049
$ws.Range("J141").Value = 3847.5
$ws.Range("K141").Value = 5315.4147
$ws.Range("L141").Value = 11542.5
$ws.Range("M141").Value = -135.4147000000003
$ws.Range("N141").Value = -21902.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1758.12
$ws.Range("I32").Value = 1600.4086
$ws.Range("J32").Value = 3853.4285
$ws.Range("K32").Value = 1600.4086
$ws.Range("L32").Value = 3853.4285
$ws.Range("M32").Value = -1313.4086
$ws.Range("N32").Value = -4427.4285

$ws.Range("H61").Value = 1973.36
$ws.Range("I61").Value = 1171.92
$ws.Range("J61").Value = 2774.8
$ws.Range("K61").Value = 1171.92
$ws.Range("L61").Value = 2774.8
$ws.Range("M61").Value = -959.9200000000001
$ws.Range("N61").Value = -3198.8

$ws.Range("H74").Value = 793.94116
$ws.Range("I74").Value = 816.5599999999999
$ws.Range("K74").Value = 816.5599999999999
$ws.Range("M74").Value = 57.44000000000005

$ws.Range("H77").Value = 793.94116
$ws.Range("I77").Value = 816.5599999999999
$ws.Range("K77").Value = 4082.8
$ws.Range("M77").Value = 285.2000000000003

$ws.Range("H132").Value = 4075
$ws.Range("I132").Value = 4546.778
$ws.Range("J132").Value = 3367.3333
$ws.Range("K132").Value = 13640.334
$ws.Range("L132").Value = 10101.9999
$ws.Range("M132").Value = -11110.334
$ws.Range("N132").Value = -15161.9999

$ws.Range("H136").Value = 1973.36
$ws.Range("I136").Value = 1171.92
$ws.Range("J136").Value = 2774.8
$ws.Range("K136").Value = 3515.76
$ws.Range("L136").Value = 8324.400000000001
$ws.Range("M136").Value = -965.7600000000002
$ws.Range("N136").Value = -13424.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 22658.143
$ws.Range("I8").Value = 269
$ws.Range("K8").Value = 269
$ws.Range("M8").Value = -129

$ws.Range("H94").Value = 500.22726
$ws.Range("I94").Value = 355.35715
$ws.Range("K94").Value = 355.35715
$ws.Range("M94").Value = 95.64285000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31437.021
$ws.Range("I31").Value = 1417.5264
$ws.Range("J31").Value = 51104.965
$ws.Range("K31").Value = 1417.5264
$ws.Range("L31").Value = 51104.965
$ws.Range("M31").Value = -1122.5264
$ws.Range("N31").Value = -51694.965

$ws.Range("H34").Value = 31437.021
$ws.Range("I34").Value = 1417.5264
$ws.Range("J34").Value = 51104.965
$ws.Range("K34").Value = 1417.5264
$ws.Range("L34").Value = 51104.965
$ws.Range("M34").Value = -1215.5264
$ws.Range("N34").Value = -51508.965

$ws.Range("H62").Value = 2399.9092
$ws.Range("I62").Value = 1849.75
$ws.Range("J62").Value = 2714.2856
$ws.Range("K62").Value = 1849.75
$ws.Range("L62").Value = 2714.2856
$ws.Range("M62").Value = -1225.75
$ws.Range("N62").Value = -3962.2856

$ws.Range("H65").Value = 2399.9092
$ws.Range("I65").Value = 1849.75
$ws.Range("J65").Value = 2714.2856
$ws.Range("K65").Value = 9248.75
$ws.Range("L65").Value = 13571.428
$ws.Range("M65").Value = -6128.75
$ws.Range("N65").Value = -19811.428

$ws.Range("H107").Value = 894.05554
$ws.Range("I107").Value = 1064.909
$ws.Range("J107").Value = 625.5714
$ws.Range("K107").Value = 1064.909
$ws.Range("L107").Value = 625.5714
$ws.Range("M107").Value = 855.0909999999999
$ws.Range("N107").Value = -4465.5714

$ws.Range("H132").Value = 30002902
$ws.Range("I132").Value = 23258422
$ws.Range("K132").Value = 69775266
$ws.Range("M132").Value = -69772736

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 578143.9
$ws.Range("J37").Value = 578143.9
$ws.Range("L37").Value = 1734431.7
$ws.Range("N37").Value = -1734655.7

$ws.Range("H116").Value = 1198
$ws.Range("I116").Value = 940.5714
$ws.Range("K116").Value = 2821.7142
$ws.Range("M116").Value = 620.2857999999997

$ws.Range("H118").Value = 2800
$ws.Range("I118").Value = 600
$ws.Range("J118").Value = 5000
$ws.Range("K118").Value = 1800
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = -557
$ws.Range("N118").Value = -17486

$ws.Range("H119").Value = 335750
$ws.Range("I119").Value = 500125
$ws.Range("J119").Value = 7000
$ws.Range("K119").Value = 1500375
$ws.Range("L119").Value = 21000
$ws.Range("M119").Value = -1495537
$ws.Range("N119").Value = -30676

$ws.Range("H120").Value = 340010
$ws.Range("I120").Value = 340010
$ws.Range("K120").Value = 1020030
$ws.Range("M120").Value = -1015192

$ws.Range("H121").Value = 7168.8887
$ws.Range("I121").Value = 6300.375
$ws.Range("J121").Value = 7863.7
$ws.Range("K121").Value = 18901.125
$ws.Range("L121").Value = 23591.1
$ws.Range("M121").Value = -17591.125
$ws.Range("N121").Value = -26211.1

$ws.Range("H131").Value = 9356.111999999999
$ws.Range("J131").Value = 9418.343000000001
$ws.Range("L131").Value = 28255.029
$ws.Range("N131").Value = -38335.029

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 19800
$ws.Range("J23").Value = 19800
$ws.Range("L23").Value = 19800
$ws.Range("N23").Value = -20246

$ws.Range("H70").Value = 54176.227
$ws.Range("I70").Value = 94534.91
$ws.Range("J70").Value = 4848.9443
$ws.Range("K70").Value = 94534.91
$ws.Range("L70").Value = 4848.9443
$ws.Range("M70").Value = -94264.91
$ws.Range("N70").Value = -5388.9443

$ws.Range("H73").Value = 54176.227
$ws.Range("I73").Value = 94534.91
$ws.Range("J73").Value = 4848.9443
$ws.Range("K73").Value = 94534.91
$ws.Range("L73").Value = 4848.9443
$ws.Range("M73").Value = -93598.91
$ws.Range("N73").Value = -6720.9443

$ws.Range("H107").Value = 631662.9
$ws.Range("I107").Value = 382
$ws.Range("J107").Value = 2525505.5
$ws.Range("K107").Value = 382
$ws.Range("L107").Value = 2525505.5
$ws.Range("M107").Value = 1538
$ws.Range("N107").Value = -2529345.5

$ws.Range("H122").Value = 778.1429000000001
$ws.Range("I122").Value = 778.1429000000001
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2334.4287
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 115.5712999999996
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 3645.6365
$ws.Range("I126").Value = 3517
$ws.Range("K126").Value = 10551
$ws.Range("M126").Value = -8081

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1133.3667
$ws.Range("I22").Value = 699.9
$ws.Range("J22").Value = 1350.1
$ws.Range("K22").Value = 699.9
$ws.Range("L22").Value = 1350.1
$ws.Range("M22").Value = -404.9
$ws.Range("N22").Value = -1940.1

$ws.Range("H27").Value = 1133.3667
$ws.Range("I27").Value = 699.9
$ws.Range("J27").Value = 1350.1
$ws.Range("K27").Value = 699.9
$ws.Range("L27").Value = 1350.1
$ws.Range("M27").Value = -592.9
$ws.Range("N27").Value = -1564.1

$ws.Range("H136").Value = 1121.225
$ws.Range("I136").Value = 944.4838999999999
$ws.Range("J136").Value = 1730
$ws.Range("K136").Value = 2833.4517
$ws.Range("L136").Value = 5190
$ws.Range("M136").Value = -283.4516999999996
$ws.Range("N136").Value = -10290
